# lifts.xlsx — trim the stale note text and clear out the long run of
# empty placeholder rows that used to pad the sheet down to row 26.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shorten the note in E3 ("Note" -> "N")
$ws.Range("E3").Value = "N"

# Move the cursor to where the next edit will happen
$null = $ws.Range("E11").Select()
